$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Cells.Item(2, 4).Value = 44447
$ws.Cells.Item(2, 10).Value = 300
$ws.Cells.Item(2, 11).Value = 900
$ws.Cells.Item(2, 12).Value = 1000
$ws.Cells.Item(2, 13).Value = 950
$ws.Cells.Item(2, 16).Value = 475

# Row 3
$ws.Cells.Item(3, 4).Value = 44525
$ws.Cells.Item(3, 10).Value = 300
$ws.Cells.Item(3, 11).Value = 1400
$ws.Cells.Item(3, 12).Value = 1500
$ws.Cells.Item(3, 13).Value = 1450
$ws.Cells.Item(3, 16).Value = 725

# Row 4
$ws.Cells.Item(4, 4).Value = 44257
$ws.Cells.Item(4, 10).Value = 500
$ws.Cells.Item(4, 11).Value = 1400
$ws.Cells.Item(4, 12).Value = 1500
$ws.Cells.Item(4, 13).Value = 1450
$ws.Cells.Item(4, 16).Value = 725

# Row 6
$ws.Cells.Item(6, 4).Value = 44202
$ws.Cells.Item(6, 10).Value = 250
$ws.Cells.Item(6, 11).Value = 1800
$ws.Cells.Item(6, 12).Value = 2000
$ws.Cells.Item(6, 13).Value = 1900
$ws.Cells.Item(6, 16).Value = 950

# Row 7
$ws.Cells.Item(7, 4).Value = 44291
$ws.Cells.Item(7, 10).Value = 250
$ws.Cells.Item(7, 11).Value = 1800
$ws.Cells.Item(7, 12).Value = 2000
$ws.Cells.Item(7, 13).Value = 1900
$ws.Cells.Item(7, 16).Value = 950

# Row 9
$ws.Cells.Item(9, 4).Value = 44572
$ws.Cells.Item(9, 11).Value = 1400
$ws.Cells.Item(9, 12).Value = 1500
$ws.Cells.Item(9, 13).Value = 1450
$ws.Cells.Item(9, 16).Value = 725

# Row 10
$ws.Cells.Item(10, 4).Value = 44392
$ws.Cells.Item(10, 10).Value = 250
$ws.Cells.Item(10, 11).Value = 1800
$ws.Cells.Item(10, 12).Value = 2000
$ws.Cells.Item(10, 13).Value = 1900
$ws.Cells.Item(10, 16).Value = 950

# Row 11
$ws.Cells.Item(11, 4).Value = 44435
$ws.Cells.Item(11, 10).Value = 300
$ws.Cells.Item(11, 11).Value = 900
$ws.Cells.Item(11, 12).Value = 1000
$ws.Cells.Item(11, 13).Value = 950
$ws.Cells.Item(11, 16).Value = 475

# Row 12
$ws.Cells.Item(12, 4).Value = 44229

# Row 13
$ws.Cells.Item(13, 4).Value = 44181
$ws.Cells.Item(13, 10).Value = 200
$ws.Cells.Item(13, 11).Value = 1000
$ws.Cells.Item(13, 12).Value = 1200
$ws.Cells.Item(13, 13).Value = 1100
$ws.Cells.Item(13, 14).Value = '$/atado'
$ws.Cells.Item(13, 16).Value = 1100
$ws.Cells.Item(13, 17).Value = 1

# Row 14
$ws.Cells.Item(14, 4).Value = 44616
$ws.Cells.Item(14, 11).Value = 1300
$ws.Cells.Item(14, 12).Value = 1500
$ws.Cells.Item(14, 13).Value = 1400
$ws.Cells.Item(14, 16).Value = 700

# Row 15
$ws.Cells.Item(15, 4).Value = 44243
$ws.Cells.Item(15, 10).Value = 250
$ws.Cells.Item(15, 11).Value = 1200
$ws.Cells.Item(15, 12).Value = 1300
$ws.Cells.Item(15, 13).Value = 1250
$ws.Cells.Item(15, 16).Value = 625

# Row 17
$ws.Cells.Item(17, 4).Value = 44302
$ws.Cells.Item(17, 10).Value = 300
$ws.Cells.Item(17, 11).Value = 900
$ws.Cells.Item(17, 12).Value = 1000
$ws.Cells.Item(17, 13).Value = 950
$ws.Cells.Item(17, 16).Value = 475

# Row 18
$ws.Cells.Item(18, 4).Value = 44795

# Row 19
$ws.Cells.Item(19, 4).Value = 44427
$ws.Cells.Item(19, 10).Value = 250
$ws.Cells.Item(19, 11).Value = 1300
$ws.Cells.Item(19, 12).Value = 1500
$ws.Cells.Item(19, 13).Value = 1400
$ws.Cells.Item(19, 16).Value = 700

# Row 20
$ws.Cells.Item(20, 4).Value = 44363
$ws.Cells.Item(20, 10).Value = 250
$ws.Cells.Item(20, 11).Value = 2500
$ws.Cells.Item(20, 12).Value = 2800
$ws.Cells.Item(20, 13).Value = 2650
$ws.Cells.Item(20, 16).Value = 1325

# Row 21
$ws.Cells.Item(21, 4).Value = 44390
$ws.Cells.Item(21, 10).Value = 250
$ws.Cells.Item(21, 11).Value = 2400
$ws.Cells.Item(21, 12).Value = 2500
$ws.Cells.Item(21, 13).Value = 2450
$ws.Cells.Item(21, 14).Value = '$/atado 1,5 a 2 kilos'
$ws.Cells.Item(21, 16).Value = 1225
$ws.Cells.Item(21, 17).Value = 2

# Row 22
$ws.Cells.Item(22, 4).Value = 44601
$ws.Cells.Item(22, 10).Value = 270
$ws.Cells.Item(22, 11).Value = 2200
$ws.Cells.Item(22, 12).Value = 2500
$ws.Cells.Item(22, 13).Value = 2350
$ws.Cells.Item(22, 16).Value = 1175

# Row 23
$ws.Cells.Item(23, 4).Value = 44365
$ws.Cells.Item(23, 10).Value = 200
$ws.Cells.Item(23, 11).Value = 1800
$ws.Cells.Item(23, 12).Value = 2000
$ws.Cells.Item(23, 13).Value = 1900
$ws.Cells.Item(23, 16).Value = 950

# Row 24
$ws.Cells.Item(24, 4).Value = 44403
$ws.Cells.Item(24, 10).Value = 250
$ws.Cells.Item(24, 11).Value = 1800
$ws.Cells.Item(24, 12).Value = 2000
$ws.Cells.Item(24, 13).Value = 1900
$ws.Cells.Item(24, 16).Value = 950

# Row 25
$ws.Cells.Item(25, 4).Value = 44385
$ws.Cells.Item(25, 11).Value = 2400
$ws.Cells.Item(25, 12).Value = 2500
$ws.Cells.Item(25, 13).Value = 2450
$ws.Cells.Item(25, 16).Value = 1225

# Row 26
$ws.Cells.Item(26, 4).Value = 44540
$ws.Cells.Item(26, 10).Value = 300

# Row 27
$ws.Cells.Item(27, 4).Value = 44726
$ws.Cells.Item(27, 11).Value = 2500
$ws.Cells.Item(27, 12).Value = 2800
$ws.Cells.Item(27, 13).Value = 2650
$ws.Cells.Item(27, 16).Value = 1325

# Row 28
$ws.Cells.Item(28, 4).Value = 44789
$ws.Cells.Item(28, 11).Value = 1400
$ws.Cells.Item(28, 12).Value = 1500
$ws.Cells.Item(28, 13).Value = 1450
$ws.Cells.Item(28, 16).Value = 725

# Row 29
$ws.Cells.Item(29, 4).Value = 44172
$ws.Cells.Item(29, 10).Value = 200
$ws.Cells.Item(29, 11).Value = 1300
$ws.Cells.Item(29, 12).Value = 1500
$ws.Cells.Item(29, 13).Value = 1400
$ws.Cells.Item(29, 16).Value = 700

# Row 30
$ws.Cells.Item(30, 4).Value = 44468
$ws.Cells.Item(30, 10).Value = 300
$ws.Cells.Item(30, 11).Value = 900
$ws.Cells.Item(30, 12).Value = 1000
$ws.Cells.Item(30, 13).Value = 950
$ws.Cells.Item(30, 16).Value = 475

# Row 31
$ws.Cells.Item(31, 4).Value = 44161
$ws.Cells.Item(31, 10).Value = 270
$ws.Cells.Item(31, 11).Value = 900
$ws.Cells.Item(31, 12).Value = 1000
$ws.Cells.Item(31, 13).Value = 950
$ws.Cells.Item(31, 16).Value = 475

# Row 32
$ws.Cells.Item(32, 4).Value = 44544
$ws.Cells.Item(32, 10).Value = 250
$ws.Cells.Item(32, 11).Value = 900
$ws.Cells.Item(32, 12).Value = 1000
$ws.Cells.Item(32, 13).Value = 950
$ws.Cells.Item(32, 16).Value = 475

# Row 33
$ws.Cells.Item(33, 4).Value = 44438
$ws.Cells.Item(33, 10).Value = 300
$ws.Cells.Item(33, 11).Value = 950
$ws.Cells.Item(33, 12).Value = 1000
$ws.Cells.Item(33, 13).Value = 975
$ws.Cells.Item(33, 16).Value = 488
